$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 4832
$ws.Range("I43").Value = 5171.143
$ws.Range("J43").Value = 4594.6
$ws.Range("K43").Value = 5171.143
$ws.Range("L43").Value = 4594.6
$ws.Range("M43").Value = -5102.143
$ws.Range("N43").Value = -4732.6

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 6278.964
$ws.Range("I112").Value = 997.25
$ws.Range("J112").Value = 6693.216
$ws.Range("K112").Value = 2991.75
$ws.Range("L112").Value = 20079.648
$ws.Range("M112").Value = -1883.75
$ws.Range("N112").Value = -22295.648

# ALC row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 2594.3845
$ws.Range("I131").Value = 2020.6364
$ws.Range("K131").Value = 6061.9092
$ws.Range("M131").Value = -1021.9092

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 9588.531999999999
$ws.Range("I137").Value = 12809.366
$ws.Range("J137").Value = 3904.7058
$ws.Range("K137").Value = 38428.098
$ws.Range("L137").Value = 11714.1174
$ws.Range("M137").Value = -35878.098
$ws.Range("N137").Value = -16814.1174

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 1114.75
$ws.Range("J5").Value = 1050
$ws.Range("L5").Value = 1050
$ws.Range("N5").Value = -1274

# ARM row 43: They've Got Legs / Steel Sabatons
$ws.Range("H43").Value = 33803.855
$ws.Range("J43").Value = 35439.332
$ws.Range("L43").Value = 35439.332
$ws.Range("N43").Value = -36065.332

# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2635.6304
$ws.Range("I61").Value = 1910.7142
$ws.Range("J61").Value = 10247.25
$ws.Range("K61").Value = 1910.7142
$ws.Range("L61").Value = 10247.25
$ws.Range("M61").Value = -1698.7142
$ws.Range("N61").Value = -10671.25

# ARM row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1001417.1
$ws.Range("I74").Value = 1001417.1
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1001417.1
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = -1000543.1
$ws.Range("M74").ClearContents()

# ARM row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1001417.1
$ws.Range("I77").Value = 1001417.1
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5007085.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5002717.5
$ws.Range("N77").Value = -5002717.5

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1322.8445
$ws.Range("I122").Value = 1239.5641
$ws.Range("K122").Value = 3718.6923
$ws.Range("M122").Value = -1268.6923

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2635.6304
$ws.Range("I136").Value = 1910.7142
$ws.Range("J136").Value = 10247.25
$ws.Range("K136").Value = 5732.142599999999
$ws.Range("L136").Value = 30741.75
$ws.Range("M136").Value = -3182.142599999999
$ws.Range("N136").Value = -35841.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 1114.75
$ws.Range("J4").Value = 1050
$ws.Range("L4").Value = 1050
$ws.Range("N4").Value = -1280

# BSM row 37: That's Some Fine Grinding / Initiate's Mortar
$ws.Range("H37").Value = 854.8889
$ws.Range("I37").Value = 727.7143
$ws.Range("K37").Value = 727.7143
$ws.Range("M37").Value = -590.7143

# BSM row 123: Archon Denied / High Durium Saw
$ws.Range("H123").Value = 65785.336
$ws.Range("I123").Value = 60000
$ws.Range("J123").Value = 88926.664
$ws.Range("K123").Value = 60000
$ws.Range("L123").Value = 88926.664
$ws.Range("N123").Value = -98726.664
$ws.Range("M123").Value = -55100

$ws = $wb.Worksheets.Item("CRP")
# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 56691.09
$ws.Range("I132").Value = 102318.164
$ws.Range("K132").Value = 306954.492
$ws.Range("M132").Value = -304424.492

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1250.8214
$ws.Range("I5").Value = 773.0909
$ws.Range("J5").Value = 3002.5
$ws.Range("K5").Value = 2319.2727
$ws.Range("L5").Value = 9007.5
$ws.Range("M5").Value = -2207.2727
$ws.Range("N5").Value = -9231.5

# CUL row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 159373.11
$ws.Range("J122").Value = 1827.6666
$ws.Range("L122").Value = 16448.9994
$ws.Range("N122").Value = -21348.9994

# CUL row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1972.037
$ws.Range("I132").Value = 1605.6666
$ws.Range("K132").Value = 14450.9994
$ws.Range("M132").Value = -11920.9994

# CUL row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1250.8214
$ws.Range("I135").Value = 773.0909
$ws.Range("J135").Value = 3002.5
$ws.Range("K135").Value = 6957.8181
$ws.Range("L135").Value = 27022.5
$ws.Range("M135").Value = -4422.8181
$ws.Range("N135").Value = -32092.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 242.5
$ws.Range("I2").Value = 292.22223
$ws.Range("J2").Value = 93.333336
$ws.Range("K2").Value = 292.22223
$ws.Range("L2").Value = 93.333336
$ws.Range("M2").Value = -179.22223
$ws.Range("N2").Value = -319.333336

# GSM row 14: All That Glitters / Copper Ear Cuffs
$ws.Range("H14").Value = 500007500
$ws.Range("I14").Value = 500007500
$ws.Range("K14").Value = 500007500
$ws.Range("M14").Value = -500007332

# GSM row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2202.2
$ws.Range("I132").Value = 2202.2
$ws.Range("K132").Value = 6606.599999999999
$ws.Range("M132").Value = -4076.599999999999

# GSM row 136: Shiny and Good / Pink Beryl
$ws.Range("H136").Value = 8572.5
$ws.Range("J136").Value = 8572.5
$ws.Range("L136").Value = 25717.5
$ws.Range("N136").Value = -30817.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 5001.722
$ws.Range("J46").Value = 10111.714
$ws.Range("L46").Value = 10111.714
$ws.Range("N46").Value = -10487.714

# LTW row 135: Dreams of Ja / Crocodileskin Leg Wraps of Scouting
$ws.Range("H135").Value = 66750
$ws.Range("J135").Value = 66750
$ws.Range("L135").Value = 66750
$ws.Range("N135").Value = -76890

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 8723
$ws.Range("I62").Value = 8999.666999999999
$ws.Range("J62").Value = 8584.666999999999
$ws.Range("K62").Value = 8999.666999999999
$ws.Range("L62").Value = 8584.666999999999
$ws.Range("M62").Value = -8375.666999999999
$ws.Range("N62").Value = -9832.666999999999

# WVR row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 8723
$ws.Range("I65").Value = 8999.666999999999
$ws.Range("J65").Value = 8584.666999999999
$ws.Range("K65").Value = 44998.335
$ws.Range("L65").Value = 42923.335
$ws.Range("M65").Value = -41878.335
$ws.Range("N65").Value = -49163.335

# WVR row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 264889.84
$ws.Range("I126").Value = 1920.5333
$ws.Range("J126").Value = 1251024.8
$ws.Range("K126").Value = 5761.5999
$ws.Range("L126").Value = 3753074.4
$ws.Range("M126").Value = -3291.5999
$ws.Range("N126").Value = -3758014.4

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 8959317
$ws.Range("I132").Value = 10033995
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 30101985
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -30099455
$ws.Range("N132").Value = -16058
